$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that was bumped by one day
# (2023-09-11 -> 2023-09-12, i.e. 45180 -> 45181) for every data row.
for ($row = 2; $row -le 224; $row++) {
    $ws.Cells.Item($row, 3).Value = 45181
}
